# Update the "取得日時" (acquisition timestamp) column (A) for data rows 2-8
# on the "ランサーズ" sheet from 2025-11-02 12:31:18 to 2025-11-02 12:42:13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-02 12:42:13"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
